# New shared-string contents introduced by this edit (oracle expression,
# javadoc tag/body, method source code, and the method-argument token) for
# the Complex#nthRoot(int) test record being appended as row 6.
$s72 = 'true ? methodResultID.stream().noneMatch(jdVar -> jdVar == null) : true;'
$s73 = '@return a List of all {@code n}-th roots of {@code this}.'
$s74 = @('/**', '     * Computes the n-th roots of this complex number.', '     * The nth roots are defined by the formula:', '     * <pre>', '     *  <code>', '     *   z<sub>k</sub> = abs<sup>1/n</sup> (cos(phi + 2&pi;k/n) + i (sin(phi + 2&pi;k/n))', '     *  </code>', '     * </pre>', '     * for <i>{@code k=0, 1, ..., n-1}</i>, where {@code abs} and {@code phi}', '     * are respectively the {@link #abs() modulus} and', '     * {@link #getArgument() argument} of this complex number.', '     * <p>', '     * If one or both parts of this complex number is NaN, a list with just', '     * one element, {@link #NaN} is returned.', '     * if neither part is NaN, but at least one part is infinite, the result', '     * is a one-element list containing {@link #INF}.', '     *', '     * @param n Degree of root.', '     * @return a List of all {@code n}-th roots of {@code this}.', '     * @throws NotPositiveException if {@code n <= 0}.', '     * @since 2.0', '     */') -join "`r`n"
$s75 = @('public List<Complex> nthRoot(int n) throws NotPositiveException {', '        if (n <= 0) {', '            throw new NotPositiveException(LocalizedFormats.CANNOT_COMPUTE_NTH_ROOT_FOR_NEGATIVE_N,', '                                           n);', '        }', '        final List<Complex> result = new ArrayList<Complex>();', '        if (isNaN) {', '            result.add(NaN);', '            return result;', '        }', '        if (isInfinite()) {', '            result.add(INF);', '            return result;', '        }', '        // nth root of abs -- faster / more accurate to use a solver here?', '        final double nthRootOfAbs = FastMath.pow(abs(), 1.0 / n);', '        // Compute nth roots of complex number with k = 0, 1, ... n-1', '        final double nthPhi = getArgument() / n;', '        final double slice = 2 * FastMath.PI / n;', '        double innerPart = nthPhi;', '        for (int k = 0; k < n ; k++) {', '            // inner part', '            final double realPart = nthRootOfAbs *  FastMath.cos(innerPart);', '            final double imaginaryPart = nthRootOfAbs *  FastMath.sin(innerPart);', '            result.add(createComplex(realPart, imaginaryPart));', '            innerPart += slice;', '        }', '        return result;', '    }') -join "`r`n"
$s76 = '[n; ; int]'

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 repeats almost all of row 5's values (same project/class/etc.), so
# duplicate row 5 first (formats, then values) - this keeps styles and the
# shared-string reuse for every unchanged column identical to a real
# "duplicate row" edit, then only the cells that actually differ are
# overwritten below.
$ws.Range("A5:V5").Copy()
$ws.Range("A6:V6").PasteSpecial(-4122)
$ws.Range("A5:V5").Copy()
$ws.Range("A6:V6").PasteSpecial(-4163)
$ws.Application.CutCopyMode = $false

# Row 5 has no content at all in columns O, Q and V - drop the blank cells
# that pasting formats created there so row 6's cell list matches row 5's.
$ws.Range("O6").Clear()
$ws.Range("Q6").Clear()
$ws.Range("V6").Clear()

# New record: id=5, a generated oracle for Complex#nthRoot(int).
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = $s72
# Leading "'" forces text entry (like typing '@return… in Excel) so the
# cell keeps the same quote-prefix style as the other javadocTag cells.
$ws.Range("G6").Value = "'" + $s73
$ws.Range("H6").Value = $s74
$ws.Range("I6").Value = $s75
$ws.Range("R6").Value = $s76

# Match the (auto-fit-to-maximum) row height used by the other wrapped rows.
$ws.Rows("6").RowHeight = 409.6

$ws.Range("U6").Select()
